$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E18").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "-"
